# Refines TODOs for search logic class
# - Adds a new time-log entry (row 37: date 3/7/2019 + new TODO note)
# - Inserts two blank/new rows (38-40) for a new "Thurs 5 am - ?" note
# - Everything previously starting at row 38 shifts down by 3 rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new rows before the old row 38, pushing existing content down.
$ws.Range("A38:A40").EntireRow.Insert()

# New content for row 37 (date + TODO note)
$ws.Range("A37").Value = 43531
$ws.Range("D37").Value = "Indi Project: consider how logic class will work - is it really a special dao?"

# New content for row 39 (follow-up note)
$ws.Range("D39").Value = "Thurs 5 am - ?"

# Match the formatting used by the other wrapped-text note cells in column D
# (the freshly inserted rows default to the time-format style copied from
# row 37's old formatting, so re-apply the plain wrap-text style here).
$ws.Range("D36").Copy()
$ws.Range("D37").PasteSpecial(-4122)
$ws.Range("D39").PasteSpecial(-4122)

# Row 38 stays blank in column D (inserting rows copies row 37's old D
# formatting into D38 by default) -- clear it out completely.
$ws.Range("D38").Clear()

# Restore the active selection to match where the author was working
$null = $ws.Range("D41").Select()
